$d = $word.ActiveDocument

$titles = @(
    "1 Image Element Title",
    "1.1 List Element Title 1",
    "1.1.1 How to take over the world",
    "1.1.2 Machining for dummies",
    "1.2 List Element Title 2",
    "1.2.1 How to boil an egg",
    "2 Image Element Title",
    "2.1 List Element Title 1",
    "2.1.1 How to take over the world",
    "2.1.2 Machining for dummies 2",
    "2.2 List Element Title 2",
    "2.2.1 How to boil an egg 2"
)

$text = [string]::Join([char]13, $titles)

$r = $d.Range(0, $d.Content.End)
$r.Text = $text

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $d.Paragraphs($i).Style = "Heading 3"
}
